$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.613.35'
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").Value = '1.923.87'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.81'
$ws.Range("E5").Value = '  +2.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2895'
$ws.Range("E8").Value = '  +1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06835'
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '105.23'
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.40'
$ws.Range("E11").Value = '  -3.79%  '
$ws.Range("D12").Value = '1.926.41'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.330'
$ws.Range("E14").Value = '  +4.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6683'
$ws.Range("E15").Value = '  +1.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '289.12'
$ws.Range("E16").Value = '  -4.64%  '
$ws.Range("D17").Value = '30.627.13'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007611'
$ws.Range("E18").Value = '  +1.61%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.580'
$ws.Range("E19").Value = '  +6.06%  '
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.95'
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("D22").Value = '2.172.88'
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.460'
$ws.Range("E24").Value = '  +2.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.509'
$ws.Range("E25").Value = '  +3.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '166.99'
$ws.Range("E26").Value = '  -0.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.13'
$ws.Range("E27").Value = '  +6.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.113'
$ws.Range("E28").Value = '  +5.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1071'
$ws.Range("E29").Value = '  -4.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.402'
$ws.Range("E30").Value = '  +3.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.171'
$ws.Range("E31").Value = '  +1.93%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("E32").Value = '  +3.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05026'
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7308'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("E35").Value = '  -0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02065'
$ws.Range("E36").Value = '  +6.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  +0.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.738'
$ws.Range("E38").Value = '  +0.37%  '
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.70'
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.049'
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8762'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4390'
$ws.Range("E43").Value = '  +6.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.928'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.68'
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.287'
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.75'
$ws.Range("E48").Value = '  +16.31%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.328'
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("E51").Value = '  +0.65%  '
